$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204699873924255
$ws.Range("B1").Value = 1.919187545776367
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.928333520889282
$ws.Range("E1").Value = 1.204297661781311
